# Update "想去人数" (interested-count) figures for the scraped 漫展 (con) entries.
# Sheet "展览" (exhibitions) and sheet "全部类型" (all types) both contain the
# same three events; bump their F-column counts to the newly scraped values.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2196   # 南宁·草莓动漫节: 2195 -> 2196
$ws1.Range("F4").Value = 1679   # 南宁·2024三月三国潮动漫节（良牙春典）: 1671 -> 1679
$ws1.Range("F5").Value = 388    # 南宁·布谷鸟动漫展4th: 387 -> 388

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2196   # 南宁·草莓动漫节: 2195 -> 2196
$ws4.Range("F6").Value = 1679   # 南宁·2024三月三国潮动漫节（良牙春典）: 1671 -> 1679
$ws4.Range("F7").Value = 388    # 南宁·布谷鸟动漫展4th: 387 -> 388
